$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.15163
$ws.Range("H2").Value = 9.454890000000001
$ws.Range("I2").Value = 0.0006291248881010851
$ws.Range("J2").Value = 0.0006291248881010851
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.94532566666667
$ws.Range("N2").Value = 77.835977
$ws.Range("O2").Value = 0.5401813355606462
$ws.Range("P2").Value = 0.5401813355606462
$ws.Range("Q2").Value = 81.77006673083667
$ws.Range("R2").Value = 735.9306005775301
$ws.Range("S2").Value = 0.0003398415222888862
$ws.Range("T2").Value = 0.0003398415222888862

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.15163
$ws.Range("H3").Value = 9.454890000000001
$ws.Range("I3").Value = 0.0006291248881010851
$ws.Range("J3").Value = 0.0006291248881010851
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("N3").Value = 55.542812
$ws.Range("O3").Value = 0.3854668692210787
$ws.Range("P3").Value = 0.3854668692210786
$ws.Range("Q3").Value = 58.35013086118668
$ws.Range("R3").Value = 525.1511777506801
$ws.Range("S3").Value = 0.0002425068009653867
$ws.Range("T3").Value = 0.0002425068009653867

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.15163
$ws.Range("H4").Value = 9.454890000000001
$ws.Range("I4").Value = 0.0006291248881010851
$ws.Range("J4").Value = 0.0006291248881010851
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.571174000000001
$ws.Range("N4").Value = 10.713522
$ws.Range("O4").Value = 0.07435179521827505
$ws.Range("P4").Value = 0.07435179521827504
$ws.Range("Q4").Value = 11.25501911362
$ws.Range("R4").Value = 101.29517202258
$ws.Range("S4").Value = 0.00004677656484681209
$ws.Range("T4").Value = 0.00004677656484681207

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4971.754394666666
$ws.Range("H5").Value = 14915.263184
$ws.Range("I5").Value = 0.9924561027819714
$ws.Range("J5").Value = 0.9924561027819713
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.94532566666667
$ws.Range("N5").Value = 77.835977
$ws.Range("O5").Value = 0.5401813355606462
$ws.Range("P5").Value = 0.5401813355606462
$ws.Range("Q5").Value = 128993.7869043079
$ws.Range("R5").Value = 1160944.082138771
$ws.Range("S5").Value = 0.5361062630860792
$ws.Range("T5").Value = 0.5361062630860792

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4971.754394666666
$ws.Range("H6").Value = 14915.263184
$ws.Range("I6").Value = 0.9924561027819714
$ws.Range("J6").Value = 0.9924561027819713
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.51427066666667
$ws.Range("N6").Value = 55.542812
$ws.Range("O6").Value = 0.3854668692210787
$ws.Range("P6").Value = 0.3854668692210786
$ws.Range("Q6").Value = 92048.40655104816
$ws.Range("R6").Value = 828435.6589594334
$ws.Range("S6").Value = 0.3825589467787196
$ws.Range("T6").Value = 0.3825589467787195

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4971.754394666666
$ws.Range("H7").Value = 14915.263184
$ws.Range("I7").Value = 0.9924561027819714
$ws.Range("J7").Value = 0.9924561027819713
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.571174000000001
$ws.Range("N7").Value = 10.713522
$ws.Range("O7").Value = 0.07435179521827505
$ws.Range("P7").Value = 0.07435179521827504
$ws.Range("Q7").Value = 17755.00002861934
$ws.Range("R7").Value = 159795.0002575741
$ws.Range("S7").Value = 0.07379089291717247
$ws.Range("T7").Value = 0.07379089291717245

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.63986933333333
$ws.Range("H8").Value = 103.919608
$ws.Range("I8").Value = 0.006914772329927541
$ws.Range("J8").Value = 0.006914772329927542
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.94532566666667
$ws.Range("N8").Value = 77.835977
$ws.Range("O8").Value = 0.5401813355606462
$ws.Range("P8").Value = 0.5401813355606462
$ws.Range("Q8").Value = 898.7426909041128
$ws.Range("R8").Value = 8088.684218137016
$ws.Range("S8").Value = 0.003735230952278061
$ws.Range("T8").Value = 0.003735230952278061

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.63986933333333
$ws.Range("H9").Value = 103.919608
$ws.Range("I9").Value = 0.006914772329927541
$ws.Range("J9").Value = 0.006914772329927542
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.51427066666667
$ws.Range("N9").Value = 55.542812
$ws.Range("O9").Value = 0.3854668692210787
$ws.Range("P9").Value = 0.3854668692210786
$ws.Range("Q9").Value = 641.3319166952996
$ws.Range("R9").Value = 5771.987250257696
$ws.Range("S9").Value = 0.002665415641393713
$ws.Range("T9").Value = 0.002665415641393713

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 34.63986933333333
$ws.Range("H10").Value = 103.919608
$ws.Range("I10").Value = 0.006914772329927541
$ws.Range("J10").Value = 0.006914772329927542
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.571174000000001
$ws.Range("N10").Value = 10.713522
$ws.Range("O10").Value = 0.07435179521827505
$ws.Range("P10").Value = 0.07435179521827504
$ws.Range("Q10").Value = 123.7050007265973
$ws.Range("R10").Value = 1113.345006539376
$ws.Range("S10").Value = 0.0005141257362557672
$ws.Range("T10").Value = 0.0005141257362557671
